$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('China',81054,46,72440,5353,1845,6,3261),
    @('Italia',59138,5560,7024,46638,3000,651,5476),
    @('Estados Unidos',32356,8149,178,31764,795,112,414),
    @('España',28603,3107,2125,24722,1785,375,1756),
    @('Alemania',24852,2488,266,24492,2,10,94),
    @('Iran',21638,1028,7635,12318,0,129,1685),
    @('Francia',16018,1559,2200,13144,1746,112,674),
    @('Corea del Sur',8897,98,2909,5884,59,2,104),
    @('Suiza',7474,611,131,7245,141,18,98),
    @('Reino Unido',5683,665,93,5309,20,48,281),
    @('Paises Bajos',4204,573,2,4023,354,43,179),
    @('Belgica',3401,586,263,3063,288,8,75),
    @('Austria',3302,310,9,3277,15,8,16),
    @('Noruega',2263,99,6,2250,28,0,7),
    @('Suecia',1931,161,16,1894,68,1,21),
    @('Portugal',1600,320,5,1581,26,2,14),
    @('Canada',1426,98,14,1392,1,1,20),
    @('Dinamarca',1395,69,1,1381,42,0,13),
    @('Australia',1353,281,46,1300,2,0,7),
    @('Malasia',1306,123,139,1157,26,2,10),
    @('Turquia',1236,289,0,1206,0,9,30),
    @('Brasil',1209,31,2,1189,18,0,18),
    @('Chequia',1120,125,6,1113,19,1,1),
    @('Japon',1086,32,235,815,57,0,36),
    @('Israel',1071,188,37,1033,18,0,1),
    @('Irlanda',906,121,5,897,29,1,4),
    @('Luxemburgo',798,128,6,784,3,0,8),
    @('Ecuador',789,257,3,772,2,7,14),
    @('Crucero',712,0,567,137,15,0,8),
    @('Pakistan',646,1,13,629,0,1,4),
    @('Polonia',634,98,13,614,3,2,1),
    @('Chile',632,95,8,623,7,0,1),
    @('Finlandia',626,103,10,615,12,0,1),
    @('Grecia',624,94,19,590,18,2,15),
    @('Tailandia',599,188,44,554,7,0,1),
    @('Islandia',568,95,5,562,1,0,1),
    @('Indonesia',514,64,29,437,0,10,48),
    @('Arabia Saudita',511,119,17,494,0,0,0),
    @('Catar',481,0,27,454,6,0,0),
    @('Singapur',455,23,144,309,14,0,2),
    @('Rumania',433,66,64,366,14,3,3),
    @('Eslovenia',414,31,0,412,12,1,2),
    @('India',396,64,24,365,0,2,7),
    @('Filipinas',380,73,15,340,1,6,25),
    @('Rusia',367,61,16,350,0,0,1),
    @('Peru',363,45,1,357,5,0,5),
    @('Barein',332,22,149,181,4,1,2),
    @('Estonia',326,20,4,322,2,0,0),
    @('Hong Kong',317,43,100,213,4,0,4),
    @('Egipto',294,0,42,242,0,0,10),
    @('Sudafrica',274,34,2,272,0,0,0),
    @('Croacia',254,48,5,248,5,0,1),
    @('Mexico',251,48,4,245,1,0,2),
    @('Libano',248,18,8,236,4,0,4),
    @('Panama',245,0,1,241,7,0,3),
    @('Irak',233,19,57,156,0,3,20),
    @('Colombia',231,35,3,226,0,2,2),
    @('Argentina',225,67,27,194,0,0,4),
    @('Serbia',222,51,2,218,4,1,2),
    @('Republica Dominicana',202,90,0,199,0,0,3),
    @('Argelia',201,62,65,119,0,2,17),
    @('Armenia',194,34,2,192,6,0,0),
    @('Kuwait',188,12,30,158,5,0,0),
    @('Bulgaria',185,22,3,179,3,0,3),
    @('Eslovaquia',185,7,7,178,2,0,0),
    @('San Marino',175,15,4,151,13,0,20),
    @('Taiwan',169,16,28,139,0,0,2),
    @('Emiratos Arabes Unidos',153,0,38,113,2,0,2),
    @('Letonia',139,15,1,138,0,0,0),
    @('Uruguay',135,25,0,135,2,0,0),
    @('Costa Rica',134,17,2,130,2,0,2),
    @('Hungria',131,28,16,109,6,2,6),
    @('Lituania',129,30,1,127,1,0,1),
    @('Republica de Macedonia',115,30,1,113,1,1,1),
    @('Islas Feroe',115,23,3,112,0,0,0),
    @('Marruecos',115,19,3,108,1,1,4),
    @('Principado de Andorra',113,25,1,111,2,1,1),
    @('Vietnam',113,19,17,96,2,0,0),
    @('Jordania',112,12,1,111,0,0,0),
    @('Republica de Chipre',95,11,3,91,3,0,1),
    @('Moldavia',94,14,1,92,3,0,1),
    @('Bosnia y Herzegovina',94,1,2,91,1,0,1),
    @('Malta',90,17,2,88,1,0,0),
    @('Albania',89,13,2,85,2,0,2),
    @('Brunei',88,5,2,86,2,0,0),
    @('Camboya',84,31,2,82,0,0,0),
    @('Sri Lanka',82,5,3,79,2,0,0),
    @('Bielorrusia',76,0,15,61,0,0,0),
    @('Tunez',75,15,1,71,7,2,3),
    @('Burkina Faso',75,11,5,66,0,1,4),
    @('Venezuela',70,0,15,55,2,0,0),
    @('Senegal',67,11,5,62,0,0,0),
    @('Nueva Zelanda',66,14,0,66,0,0,0),
    @('Azerbaiyan',65,12,11,53,0,0,1),
    @('Ucrania',63,16,1,59,0,0,3),
    @('Kazajistan',59,5,0,59,0,0,0),
    @('Estado de Palestina',59,6,17,42,0,0,0),
    @('Guadalupe',56,0,0,38,4,0,1),
    @('Oman',55,3,17,38,0,0,0),
    @('Georgia',54,5,1,53,1,0,0),
    @('Trinidad yTobago',50,1,0,50,0,0,0),
    @('Reunion',47,0,1,46,0,0,0),
    @('Uzbekistan',43,2,0,43,0,0,0),
    @('Camerun',40,0,2,38,0,0,0),
    @('Afganistan',40,16,1,38,0,1,1),
    @('Liechtenstein',37,0,0,37,0,0,0),
    @('Martinica',37,0,0,36,7,0,1),
    @('Cuba',35,14,0,34,0,0,1),
    @('Consejo Danes para los Refugiados',30,7,0,29,0,0,1),
    @('Nigeria',30,8,2,28,0,0,0),
    @('Guam',27,12,0,26,0,1,1),
    @('Banglades',27,3,3,22,0,0,2),
    @('Honduras',26,2,0,26,0,0,0),
    @('Bolivia',24,5,0,24,0,0,0),
    @('Mauricio',24,10,0,22,1,1,2),
    @('Monaco',23,5,1,22,0,0,0),
    @('Puerto Rico',23,2,0,22,0,0,1),
    @('Ghana',23,2,0,21,0,0,1),
    @('Paraguay',22,0,0,21,1,0,1),
    @('Macao',22,3,10,12,0,0,0),
    @('Montenegro',21,5,0,21,0,0,0),
    @('Jamaica',19,0,2,16,0,0,1),
    @('Guayana Francesa',18,0,0,18,0,0,0),
    @('Guyana',18,0,0,17,0,0,1),
    @('Ruanda',17,0,0,17,0,0,0),
    @('Guatemala',17,0,0,16,0,0,1),
    @('Togo',16,0,0,16,0,0,0),
    @('Polinesia Francesa',15,0,0,15,0,0,0),
    @('Kenia',15,8,0,15,0,0,0),
    @('Gibraltar',15,5,2,13,0,0,0),
    @('Kirguistan',14,0,0,14,0,0,0),
    @('Barbados',14,0,0,14,0,0,0),
    @('Costa de Marfil',14,0,1,13,0,0,0),
    @('Maldivas',13,0,3,10,0,0,0),
    @('Tanzania',12,6,0,12,0,0,0),
    @('Mayotte',11,0,0,11,0,0,0),
    @('Etiopia',11,2,0,11,0,0,0),
    @('Mongolia',10,0,0,10,0,0,0),
    @('Aruba',8,3,1,7,0,0,0),
    @('Seychelles',7,0,0,7,0,0,0),
    @('Islas Virgenes de los Estados Unidos',6,0,0,6,0,0,0),
    @('Guinea Ecuatorial',6,0,0,6,0,0,0),
    @('Surinam',5,0,0,5,0,0,0),
    @('Isla de Man',5,3,0,5,0,0,0),
    @('San Martin (Parte Francesa)',5,0,0,5,0,0,0),
    @('Gabon',5,0,0,4,0,0,1),
    @('Suazilandia',4,3,0,4,0,0,0),
    @('Nueva Caledonia',4,0,0,4,0,0,0),
    @('Bahamas',4,0,0,4,0,0,0),
    @('Namibia',3,0,0,3,0,0,0),
    @('Cabo Verde',3,0,0,3,0,0,0),
    @('San Bartolome',3,0,0,3,0,0,0),
    @('Congo',3,0,0,3,0,0,0),
    @('Liberia',3,0,0,3,0,0,0),
    @('Madagascar',3,0,0,3,0,0,0),
    @('Zimbabue',3,0,0,3,0,0,0),
    @('Republica de Africa Central',3,0,0,3,0,0,0),
    @('El Salvador',3,0,0,3,0,0,0),
    @('Zambia',3,1,0,3,0,0,0),
    @('Curazao',3,0,0,2,0,0,1),
    @('Islas Caimanes',3,0,0,2,0,0,1),
    @('Angola',2,0,0,2,0,0,0),
    @('Butan',2,0,0,2,0,0,0),
    @('Groenlandia',2,0,0,2,0,0,0),
    @('Fiyi',2,0,0,2,0,0,0),
    @('Santa Lucia',2,0,0,2,0,0,0),
    @('Guinea',2,0,0,2,0,0,0),
    @('Bermudas',2,0,0,2,0,0,0),
    @('Benin',2,0,0,2,0,0,0),
    @('Nicaragua',2,0,0,2,0,0,0),
    @('Mauritania',2,0,0,2,0,0,0),
    @('Niger',2,1,0,2,0,0,0),
    @('Haiti',2,0,0,2,0,0,0),
    @('Sudan',2,0,0,1,0,0,1),
    @('Siria',1,1,0,1,0,0,0),
    @('Granada',1,1,0,1,0,0,0),
    @('San Martin (Parte Holandesa)',1,0,0,1,0,0,0),
    @('Montserrat',1,0,0,1,0,0,0),
    @('Gambia',1,0,0,1,0,0,0),
    @('Santa Sede',1,0,0,1,0,0,0),
    @('Papua Nueva Guinea',1,0,0,1,0,0,0),
    @('San Vicente y las Granadinas',1,0,0,1,0,0,0),
    @('Timor Oriental',1,0,0,1,0,0,0),
    @('Eritrea',1,0,0,1,0,0,0),
    @('Uganda',1,0,0,1,0,0,0),
    @('Somalia',1,0,0,1,0,0,0),
    @('Republica de Yibuti',1,0,0,1,0,0,0),
    @('Mozambique',1,1,0,1,0,0,0),
    @('Antigua y Barbuda',1,0,0,1,0,0,0),
    @('Republica del Chad',1,0,0,1,0,0,0),
    @('Nepal',1,0,1,0,0,0,0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 4 + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = [double]$rowData[$c]
    }
}

$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 21:16"

Write-Output "done"
